$d = $word.ActiveDocument

# --- Change 1: rewrite the opening paragraph ---
$d.Content.Find.Execute(
    "Uma Rádio Comunitária sem fins comercais necessita de um software para que possa realizar as operações normais e rotineiras de um Player de Áudio comum, mais ainda possibilitando uma maior autonomia.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Uma Rádio necessita de um software para contemplar operações rotineiras, como, por exemplo, reproduzir arquivos de áudio, automatizar eventos quando não houver um Operador de Áudio.",
    2) | Out-Null

# --- Change 2: split the existing blank paragraph into two blank
#     paragraphs by inserting a paragraph mark right before the "O
#     Sistema" text (matching only a short prefix keeps the run's
#     character formatting intact) ---
$d.Content.Find.Execute(
    "O Sistema deve contemplar",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "^pO Sistema deve contemplar",
    2) | Out-Null

# --- Change 3: now rewrite the "O Sistema" paragraph text (separate
#     call so the run keeps its <w:rPr>) ---
$d.Content.Find.Execute(
    "O Sistema deve contemplar toda a rotina de uma rádio, reproduzir músicas diversas, vinhetas, spot, programações gravadas, mantendo a Playlist e eventos automatizados em banco de dados. Precisa ser organizado por Módulos, sendo o principal, um Player de Áudio comum, contendo todas as operações básicas de um player (Play, Stop, Playlist, execução manual, sequencial, aleatória) sendo capaz de executar músicas em formato MP3.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "O Sistema deve reproduzir músicas diversas, vinhetas, spot, programações gravadas, mantendo a Playlist e eventos automatizados em banco de dados. Precisa ser organizado por Módulos, sendo o principal, um Player de Áudio comum, contendo todas as operações básicas, sendo capaz de executar músicas em formato MP3.",
    2) | Out-Null

# --- Change 4: rewrite the "Há necessidade" paragraph ---
$d.Content.Find.Execute(
    "Há necessidade de haver um Módulo para a inserção das Vinhetas e Spots, para o gerenciamento dos eventos automatizados, outro para as programações gravadas, sendo possível reproduzi-las de maneira eficiente, redirecionando a Playlist ao término.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Há necessidade de Módulos bem definidos, para a inserção das Vinhetas e Spots, para o gerenciamento dos eventos automatizados, para o gerenciamento das programações gravadas, sendo possível reproduzi-los de maneira eficiente, redirecionando a Playlist ao término.",
    2) | Out-Null

# --- Change 5: merge the final paragraph's two runs into one, dropping
#     the run that follows the _GoBack bookmark ---
$bm = $d.Bookmarks("_GoBack")
$para = $bm.Range.Paragraphs(1)
$afterBookmark = $d.Range($bm.End, $para.Range.End - 1)
$afterBookmark.Text = ""

$d.Content.Find.Execute(
    "O sistema",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "O sistema deve ser capaz de funcionar em outras plataformas, como Windows, Linux, dando maior liberdade para a escolha da plataforma, deve funcionar 24 horas por dia, caso ele seja fechado, permanecer em execução em segundo plano, nunca sendo totalmente encerrado.",
    2) | Out-Null
